# Ticket 48 - add an "Immaterial" sheet whose template exercises a
# countdown-by-step jt:for loop, and make sure that blank-but-styled
# trailer cells (row 3) are present so that shifting logic has to take
# their non-default CellStyle into account.

$wb = $excel.ActiveWorkbook

# Add the new worksheet as the last tab in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Immaterial"

# Header row.
$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "square"

# Template row driving the for-loop.
$ws.Range("A2").Value = '<jt:for start="10" end="0" step="-2" var="n">${n}'
$ws.Range("B2").Value = '${n * n}</jt:for>'

# Style the header row: bold 9pt Arial on a yellow fill.
$headerFont = $ws.Range("A1:B1").Font
$headerFont.Bold = $true
$headerFont.Name = "Arial"
$headerFont.Size = 9
$ws.Range("A1:B1").Interior.Color = 65535

# Row 3 is intentionally left blank, but still carries a (red) fill so
# that it counts as "styled" content for shifting purposes.
$ws.Range("A3:B3").Interior.Color = 255

# Size the columns to fit their (long) template contents.
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
